# Auto-generated edit script: updates cryptos list values per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.694.02"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.615.69"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.992"
$ws.Range("E4").Value = "  -0.63%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.67"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.991"
$ws.Range("E7").Value = "  -0.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.22"
$ws.Range("E8").Value = "  +10.01%  "
$ws.Range("E9").Value = "  +3.12%  "
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.845.55"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.606.52"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.570"
$ws.Range("E14").Value = "  +6.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.92"
$ws.Range("E15").Value = "  +6.12%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "9.04"
$ws.Range("E16").Value = "  +18.39%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.677.34"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.19"
$ws.Range("E18").Value = "  +1.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.75"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0712"
$ws.Range("E20").Value = "  +3.33%  "
$ws.Range("E21").Value = "  -0.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.11"
$ws.Range("E22").Value = "  +3.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.70"
$ws.Range("E23").Value = "  +6.19%  "
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.42"
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.67"
$ws.Range("E26").Value = "  +2.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.111"
$ws.Range("E27").Value = "  +2.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.60"
$ws.Range("E28").Value = "  +3.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.993"
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0488"
$ws.Range("E30").Value = "  +3.55%  "
$ws.Range("E31").Value = "  +3.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.33"
$ws.Range("E32").Value = "  +3.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.22"
$ws.Range("E33").Value = "  +4.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.425.66"
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("E35").Value = "  +7.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.05"
$ws.Range("E36").Value = "  +1.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.86"
$ws.Range("E37").Value = "  +1.87%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0171"
$ws.Range("E38").Value = "  +3.29%  "
$ws.Range("B39").Value = "HuobiToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.28"
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.558"
$ws.Range("E40").Value = "  +4.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0506"
$ws.Range("E41").Value = "  +3.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.832"
$ws.Range("E42").Value = "  +4.69%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "69.86"
$ws.Range("E44").Value = "  +6.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "53.92"
$ws.Range("E45").Value = "  +2.31%  "
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("E47").Value = "  +18.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.45"
$ws.Range("E48").Value = "  +3.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.754.77"
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "87.95"
$ws.Range("E50").Value = "  +1.51%  "
$ws.Range("E51").Value = "  -1.72%  "
